$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Derek Bays" -> "Derek J Bays" in author list; fill in B2 affiliations; M2 17->19
$ws.Cells.Item(2, 1).Value = 'Amy P. Hsu, Agnieszka Korzeniowska, Cynthia C Aguilar, Jiande Gu, Eric Karlins, Andrew J. Oler, Gang Chen, Glennys V. Reynoso, Joie Davis, Alexandria Laurel Chaput, Tao Peng, Ling Sun, Justin Lack, Derek J Bays, Ethan R. Stewart, Sarah Waldman, Daniel A. Powell, Fariba M. Donovan, Jigar V. Desai, Nima Pouladi, Debra A. Long Priel, Daisuke Yamanaka, Sergio D. Rosenzweig, Julie E. Niemela, Jennifer Stoddard, Alexandra F. Freeman, Christa S. Zerbe, Douglas B. Kuhns, Yves A. Lussier, Kenneth N. Olivier, Richard C. Boucher, Heather D. Hickman, Jeffrey A. Frelinger, Joshua Fierer, Lisa F. Shubitz, Thomas L. Leto, George R. Thompson, John N. Galgiani, Michail S Lionakis, Steven M. Holland'
$ws.Cells.Item(2, 2).Value = 'Department of Cell Biology and Molecular Genetics, University of Maryland, College Park, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; ; ; ; Marsico Lung Institute and Cystic Fibrosis Research Center, University of North Carolina at Chapel Hill, Chapel Hill, North Carolina, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Valley Fever Center for Excellence, University of Arizona College of Medicine–Tucson, Tucson, Arizona, USA.; Valley Fever Center for Excellence, University of Arizona College of Medicine–Tucson, Tucson, Arizona, USA.; Department of Respiratory and Critical Care Medicine, Laboratory of Pulmonary Immunology and Inflammation, West China Hospital, Sichuan University, Chengdu, Sichuan Province, China.; Marsico Lung Institute and Cystic Fibrosis Research Center, University of North Carolina at Chapel Hill, Chapel Hill, North Carolina, USA.; Advanced Biomedical Computational Science, Frederick National Laboratory for Cancer Research, Leidos Biomedical Research, Inc., Frederick, Maryland, USA.; NIAID Collaborative Bioinformatics Resource, NIAID, NIH, Bethesda, Maryland, USA.; Department of Internal Medicine, Division of Infectious Diseases, UC Davis Health, Sacramento, California, USA.; Department of Internal Medicine, Division of Infectious Diseases, UC Davis Health, Sacramento, California, USA.; Department of Internal Medicine, Division of Infectious Diseases, UC Davis Health, Sacramento, California, USA.; Department of Immunobiology, University of Arizona, Tucson, Arizona, USA.; Valley Fever Center for Excellence, University of Arizona College of Medicine–Tucson, Tucson, Arizona, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Center for Biomedical Informatics and Biostatistics and; The Center for Applied Genetics and Genomic Medicine, Department of Medicine, University of Arizona, Tucson, Arizona, USA.; Neutrophil Monitoring Laboratory, Applied/Developmental Research Directorate, Leidos Biomedical Research, Inc, Frederick National Laboratory for Cancer Research, Frederick, Maryland, USA.; Laboratory for Immunopharmacology of Microbial Products, School of Pharmacy, Tokyo University of Pharmacy and Life Sciences, Hachioji, Tokyo, Japan.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; ; ; ; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Neutrophil Monitoring Laboratory, Applied/Developmental Research Directorate, Leidos Biomedical Research, Inc, Frederick National Laboratory for Cancer Research, Frederick, Maryland, USA.; Center for Biomedical Informatics and Biostatistics and; The Center for Applied Genetics and Genomic Medicine, Department of Medicine, University of Arizona, Tucson, Arizona, USA.; Laboratory of Chronic Airway Infection, Pulmonary Branch, National Heart, Lung, and Blood Institute, NIH, Bethesda, Maryland, USA.; Marsico Lung Institute and Cystic Fibrosis Research Center, University of North Carolina at Chapel Hill, Chapel Hill, North Carolina, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Department of Immunobiology, University of Arizona, Tucson, Arizona, USA.; Division of Infectious Diseases, Departments of Pathology and Medicine, School of Medicine, University of California San Diego, La Jolla, California, USA.; VA HealthCare San Diego, San Diego, California, USA.; Valley Fever Center for Excellence, University of Arizona College of Medicine–Tucson, Tucson, Arizona, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Department of Medical Microbiology and Immunology, University of California Davis, Davis, California, USA.; Valley Fever Center for Excellence, University of Arizona College of Medicine–Tucson, Tucson, Arizona, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.'
$ws.Cells.Item(2, 13).NumberFormat = "@"
$ws.Cells.Item(2, 13).Value = '19'

# Rows 4 and 5: article data content is swapped between the two rows; M (cited_by_count) gets fresh values
$ws.Cells.Item(4, 1).Value = 'M. Alejandra Mandel, Sinem Beyhan, Mark Voorhies, Lisa F. Shubitz, John N. Galgiani, Marc J. Orbach, Anita Sil'
$ws.Cells.Item(4, 2).Value = 'School of Plant Sciences, University of Arizona, Tucson, Arizona, United States of America.; Department of Microbiology and Immunology, University of California San Francisco, San Francisco, California, United States of America; Department of Microbiology and Immunology, University of California San Francisco, San Francisco, California, United States of America; Valley Fever Center for Excellence, University of Arizona, Tucson, Arizona, United States of America; Valley Fever Center for Excellence, University of Arizona, Tucson, Arizona, United States of America; School of Plant Sciences, University of Arizona, Tucson, Arizona, United States of America; Department of Microbiology and Immunology, University of California San Francisco, San Francisco, California, United States of America'
$ws.Cells.Item(4, 3).Value = 'https://openalex.org/W4225503377'
$ws.Cells.Item(4, 4).Value = 'The WOPR family protein Ryp1 is a key regulator of gene expression, development, and virulence in the thermally dimorphic fungal pathogen Coccidioides posadasii'
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '2022-04-06'
$ws.Cells.Item(4, 6).Value = 'PLOS Pathogens'
$ws.Cells.Item(4, 7).Value = 'Public Library of Science'
$ws.Cells.Item(4, 8).Value = 'https://doi.org/10.1371/journal.ppat.1009832'
$ws.Cells.Item(4, 9).Value = 'cc-by'
$ws.Cells.Item(4, 10).Value = 'publishedVersion'
$ws.Cells.Item(4, 11).Value = 'gold'
$ws.Cells.Item(4, 15).Value = 'https://pubmed.ncbi.nlm.nih.gov/35385558'
$ws.Cells.Item(4, 16).Value = 'https://doi.org/10.1371/journal.ppat.1009832'
$ws.Cells.Item(5, 1).Value = 'George R. Thompson, Neil M. Ampel, Janis E. Blair, Fariba M. Donovan, Joshua Fierer, John N. Galgiani, Arash Heidari, Royce H. Johnson, Stanley A. Shatsky, Christopher M. Uchiyama, David A. Stevens'
$ws.Cells.Item(5, 2).Value = 'Department of Internal Medicine, Division of Infectious Diseases and the Department of Medical Microbiology and Immunology, University of California-Davis Medical Center;  Sacramento, CA  USA.; University of California - Davis Center for Valley Fever, Sacramento, CA  USA.; Division of Infectious Diseases, Mayo Clinic in Arizona,  Phoenix, AZ,  USA.; University of Arizona College of Medicine,  Tucson, AZ,  USA.; Division of Infectious Diseases, Mayo Clinic in Arizona,  Phoenix, AZ,  USA.; University of Arizona College of Medicine,  Tucson, AZ,  USA.; Division of Infectious Disease, Department of Medicine, University of California San Diego School of Medicine, La Jolla, CA,  USA.; Infectious Diseases Section, VA Healthcare San Diego,  San Diego, CA, USA.; University of Arizona College of Medicine,  Tucson, AZ,  USA.; Division of Infectious Diseases, Department of Medicine, David Geffen School of Medicine at UCLA, Kern Medical,  Bakersfield, CA, USA.; Division of Infectious Diseases, Department of Medicine, David Geffen School of Medicine at UCLA, Kern Medical,  Bakersfield, CA, USA.; Good Samaritan Hospital, San Jose, CA,  USA.; Department of Neurosurgery, Scripps Clinic and Scripps Green Hospital, La Jolla, CA, USA.; California Institute for Medical Research, San Jose, CA, USA.; Division of Infectious Diseases and Geographic Medicine, Department of Medicine, Stanford University School of Medicine, Stanford, CA, USA.'
$ws.Cells.Item(5, 3).Value = 'https://openalex.org/W4283122802'
$ws.Cells.Item(5, 4).Value = 'Controversies in the Management of Central Nervous System Coccidioidomycosis'
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '2022-06-19'
$ws.Cells.Item(5, 6).Value = 'Clinical Infectious Diseases'
$ws.Cells.Item(5, 7).Value = 'Oxford University Press'
$ws.Cells.Item(5, 8).Value = 'https://doi.org/10.1093/cid/ciac478'
$ws.Cells.Item(5, 9).Value = 'N/A'
$ws.Cells.Item(5, 10).Value = 'N/A'
$ws.Cells.Item(5, 11).Value = 'closed'
$ws.Cells.Item(5, 15).Value = 'https://pubmed.ncbi.nlm.nih.gov/35717645'
$ws.Cells.Item(5, 16).Value = 'https://doi.org/10.1093/cid/ciac478'
$ws.Cells.Item(4, 13).NumberFormat = "@"
$ws.Cells.Item(4, 13).Value = '9'
$ws.Cells.Item(5, 13).NumberFormat = "@"
$ws.Cells.Item(5, 13).Value = '8'

# Rows 6 and 7: article data content is swapped between the two rows; M (cited_by_count) gets fresh values
$ws.Cells.Item(6, 1).Value = 'John N. Galgiani, Lisa F. Shubitz, Marc J. Orbach, M. Alejandra Mandel, Daniel A. Powell, Bruce S. Klein, Edward J. Robb, Mana Ohkura, Devin J. Seka, Thomas Tomasiak, Thomas P. Monath'
$ws.Cells.Item(6, 2).Value = 'Valley Fever Center for Excellence, University of Arizona College of Medicine-Tucson, Tucson, AZ 85724, USA; Valley Fever Center for Excellence, University of Arizona College of Medicine-Tucson, Tucson, AZ 85724, USA; Valley Fever Center for Excellence, University of Arizona College of Medicine-Tucson, Tucson, AZ 85724, USA; Valley Fever Center for Excellence, University of Arizona College of Medicine-Tucson, Tucson, AZ 85724, USA; Valley Fever Center for Excellence, University of Arizona College of Medicine-Tucson, Tucson, AZ 85724, USA; Department of Pediatrics, University of Wisconsin School of Medicine and Public Health, University of Wisconsin-Madison, Madison, WI 53706, USA; Anivive Lifesciences, Long Beach, CA 90807, USA; Department of Botany and Plant Pathology, College of Agricultural Sciences, Oregon State University, Corvallis, OR 97331, USA; Department of Chemistry and Biochemistry, College of Science, University of Arizona, Tucson, AZ 85721, USA; Bio5 Institute, University of Arizona, Tucson, AZ 85721, USA; Crozet Biopharma LLC, Lexington, MA 02420, USA'
$ws.Cells.Item(6, 3).Value = 'https://openalex.org/W4291001696'
$ws.Cells.Item(6, 4).Value = 'Vaccines to Prevent Coccidioidomycosis: A Gene-Deletion Mutant of Coccidioides Posadasii as a Viable Candidate for Human Trials'
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '2022-08-10'
$ws.Cells.Item(6, 6).Value = 'Journal of Fungi'
$ws.Cells.Item(6, 7).Value = 'Multidisciplinary Digital Publishing Institute'
$ws.Cells.Item(6, 8).Value = 'https://doi.org/10.3390/jof8080838'
$ws.Cells.Item(6, 9).Value = 'cc-by'
$ws.Cells.Item(6, 10).Value = 'publishedVersion'
$ws.Cells.Item(6, 11).Value = 'gold'
$ws.Cells.Item(6, 15).Value = 'https://pubmed.ncbi.nlm.nih.gov/36012826'
$ws.Cells.Item(6, 16).Value = 'https://doi.org/10.3390/jof8080838'
$ws.Cells.Item(7, 1).Value = 'Daniel A. Powell, Amy P. Hsu, Christine D. Butkiewicz, Hien T. Trinh, Jeffrey A. Frelinger, Steven M. Holland, John N. Galgiani, Lisa F. Shubitz'
$ws.Cells.Item(7, 2).Value = 'Valley Fever Center for Excellence, University of Arizona, Tucson, AZ, United States; Laboratory of Clinical and Infectious Diseases, National Institutes of Allergy and Infectious Disease, Bethesda, MD, United States; Valley Fever Center for Excellence, University of Arizona, Tucson, AZ, United States; Valley Fever Center for Excellence, University of Arizona, Tucson, AZ, United States; Valley Fever Center for Excellence, University of Arizona, Tucson, AZ, United States; Laboratory of Clinical and Infectious Diseases, National Institutes of Allergy and Infectious Disease, Bethesda, MD, United States; Valley Fever Center for Excellence, University of Arizona, Tucson, AZ, United States; Valley Fever Center for Excellence, University of Arizona, Tucson, AZ, United States'
$ws.Cells.Item(7, 3).Value = 'https://openalex.org/W4205456491'
$ws.Cells.Item(7, 4).Value = 'Vaccine Protection of Mice With Primary Immunodeficiencies Against Disseminated Coccidioidomycosis'
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '2022-01-07'
$ws.Cells.Item(7, 6).Value = 'Frontiers in Cellular and Infection Microbiology'
$ws.Cells.Item(7, 7).Value = 'Frontiers Media'
$ws.Cells.Item(7, 8).Value = 'https://doi.org/10.3389/fcimb.2021.790488'
$ws.Cells.Item(7, 9).Value = 'cc-by'
$ws.Cells.Item(7, 10).Value = 'publishedVersion'
$ws.Cells.Item(7, 11).Value = 'gold'
$ws.Cells.Item(7, 15).Value = 'https://pubmed.ncbi.nlm.nih.gov/35071044'
$ws.Cells.Item(7, 16).Value = 'https://doi.org/10.3389/fcimb.2021.790488'
$ws.Cells.Item(6, 13).NumberFormat = "@"
$ws.Cells.Item(6, 13).Value = '6'
$ws.Cells.Item(7, 13).NumberFormat = "@"
$ws.Cells.Item(7, 13).Value = '5'

# Row 10: "Derek Bays" -> "Derek J Bays"
$ws.Cells.Item(10, 1).Value = 'Joey Shemuel, Derek J Bays, George R. Thompson, Susan E. Reef, Linda Snyder, Alana Freifeld, Milt Huppert, David Salkin, Machelle Wilson, John N. Galgiani'

# Row 11: M11 1 -> 2
$ws.Cells.Item(11, 13).NumberFormat = "@"
$ws.Cells.Item(11, 13).Value = '2'
